# Switched labels for 8A and 8B so that they are in chronological order
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 86-89 currently hold "8A" in column A -> should become "8B"
# Rows 90-93 currently hold "8B" in column A -> should become "8A"
for ($r = 86; $r -le 89; $r++) {
    $ws.Cells.Item($r, 1).Value = "8B"
}
for ($r = 90; $r -le 93; $r++) {
    $ws.Cells.Item($r, 1).Value = "8A"
}

# Update the view to reflect where the user ended up scrolled/selected after the edit
$ws.Range("A94").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 67
$win.ScrollColumn = 1
